$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(47, 1).Value = 'Create Country'
$ws.Cells.Item(47, 2).Value = 'PASSED'
$ws.Cells.Item(47, 3).Value = 'chrome'
$ws.Cells.Item(47, 4).Value = '13.12.22'
$ws.Cells.Item(48, 1).Value = 'Add Document Types'
$ws.Cells.Item(48, 2).Value = 'PASSED'
$ws.Cells.Item(48, 3).Value = 'chrome'
$ws.Cells.Item(48, 4).Value = '16.12.22'
$ws.Cells.Item(49, 1).Value = 'Add Document Types'
$ws.Cells.Item(49, 2).Value = 'PASSED'
$ws.Cells.Item(49, 3).Value = 'chrome'
$ws.Cells.Item(49, 4).Value = '17.12.22'
$ws.Cells.Item(50, 1).Value = 'Edit Document Types'
$ws.Cells.Item(50, 2).Value = 'FAILED'
$ws.Cells.Item(50, 3).Value = 'chrome'
$ws.Cells.Item(50, 4).Value = '17.12.22'
$ws.Cells.Item(51, 1).Value = 'Edit Document Types'
$ws.Cells.Item(51, 2).Value = 'FAILED'
$ws.Cells.Item(51, 3).Value = 'chrome'
$ws.Cells.Item(51, 4).Value = '17.12.22'
$ws.Cells.Item(52, 1).Value = 'Edit Document Types'
$ws.Cells.Item(52, 2).Value = 'FAILED'
$ws.Cells.Item(52, 3).Value = 'chrome'
$ws.Cells.Item(52, 4).Value = '17.12.22'
$ws.Cells.Item(53, 1).Value = 'Edit Document Types'
$ws.Cells.Item(53, 2).Value = 'FAILED'
$ws.Cells.Item(53, 3).Value = 'chrome'
$ws.Cells.Item(53, 4).Value = '17.12.22'
$ws.Cells.Item(54, 1).Value = 'Edit Document Types'
$ws.Cells.Item(54, 2).Value = 'FAILED'
$ws.Cells.Item(54, 3).Value = 'chrome'
$ws.Cells.Item(54, 4).Value = '18.12.22'
$ws.Cells.Item(55, 1).Value = 'Edit Document Types'
$ws.Cells.Item(55, 2).Value = 'FAILED'
$ws.Cells.Item(55, 3).Value = 'chrome'
$ws.Cells.Item(55, 4).Value = '18.12.22'
$ws.Cells.Item(56, 1).Value = 'Add Document Types'
$ws.Cells.Item(56, 2).Value = 'PASSED'
$ws.Cells.Item(56, 3).Value = 'chrome'
$ws.Cells.Item(56, 4).Value = '18.12.22'
$ws.Cells.Item(57, 1).Value = 'Edit Document Types'
$ws.Cells.Item(57, 2).Value = 'FAILED'
$ws.Cells.Item(57, 3).Value = 'chrome'
$ws.Cells.Item(57, 4).Value = '18.12.22'
$ws.Cells.Item(58, 1).Value = 'Add Document Types'
$ws.Cells.Item(58, 2).Value = 'PASSED'
$ws.Cells.Item(58, 3).Value = 'chrome'
$ws.Cells.Item(58, 4).Value = '18.12.22'
$ws.Cells.Item(59, 1).Value = 'Edit Document Types'
$ws.Cells.Item(59, 2).Value = 'FAILED'
$ws.Cells.Item(59, 3).Value = 'chrome'
$ws.Cells.Item(59, 4).Value = '18.12.22'
$ws.Cells.Item(60, 1).Value = 'Edit Document Types'
$ws.Cells.Item(60, 2).Value = 'PASSED'
$ws.Cells.Item(60, 3).Value = 'chrome'
$ws.Cells.Item(60, 4).Value = '18.12.22'
$ws.Cells.Item(61, 1).Value = 'Add Document Types'
$ws.Cells.Item(61, 2).Value = 'PASSED'
$ws.Cells.Item(61, 3).Value = 'chrome'
$ws.Cells.Item(61, 4).Value = '18.12.22'
$ws.Cells.Item(62, 1).Value = 'Edit Document Types'
$ws.Cells.Item(62, 2).Value = 'FAILED'
$ws.Cells.Item(62, 3).Value = 'chrome'
$ws.Cells.Item(62, 4).Value = '18.12.22'
$ws.Cells.Item(63, 1).Value = 'Add Document Types'
$ws.Cells.Item(63, 2).Value = 'PASSED'
$ws.Cells.Item(63, 3).Value = 'chrome'
$ws.Cells.Item(63, 4).Value = '19.12.22'
$ws.Cells.Item(64, 1).Value = 'Edit Document Types'
$ws.Cells.Item(64, 2).Value = 'FAILED'
$ws.Cells.Item(64, 3).Value = 'chrome'
$ws.Cells.Item(64, 4).Value = '19.12.22'
$ws.Cells.Item(65, 1).Value = 'Delete Document Types'
$ws.Cells.Item(65, 2).Value = 'FAILED'
$ws.Cells.Item(65, 3).Value = 'chrome'
$ws.Cells.Item(65, 4).Value = '19.12.22'
$ws.Cells.Item(66, 1).Value = 'Edit Document Types'
$ws.Cells.Item(66, 2).Value = 'FAILED'
$ws.Cells.Item(66, 3).Value = 'chrome'
$ws.Cells.Item(66, 4).Value = '19.12.22'
$ws.Cells.Item(67, 1).Value = 'Edit Document Types'
$ws.Cells.Item(67, 2).Value = 'FAILED'
$ws.Cells.Item(67, 3).Value = 'chrome'
$ws.Cells.Item(67, 4).Value = '19.12.22'
$ws.Cells.Item(68, 1).Value = 'Edit Document Types'
$ws.Cells.Item(68, 2).Value = 'PASSED'
$ws.Cells.Item(68, 3).Value = 'chrome'
$ws.Cells.Item(68, 4).Value = '19.12.22'
$ws.Cells.Item(69, 1).Value = 'Delete Document Types'
$ws.Cells.Item(69, 2).Value = 'PASSED'
$ws.Cells.Item(69, 3).Value = 'chrome'
$ws.Cells.Item(69, 4).Value = '19.12.22'
$ws.Cells.Item(70, 1).Value = 'Add Document Types'
$ws.Cells.Item(70, 2).Value = 'PASSED'
$ws.Cells.Item(70, 3).Value = 'chrome'
$ws.Cells.Item(70, 4).Value = '19.12.22'
$ws.Cells.Item(71, 1).Value = 'Edit Document Types'
$ws.Cells.Item(71, 2).Value = 'FAILED'
$ws.Cells.Item(71, 3).Value = 'chrome'
$ws.Cells.Item(71, 4).Value = '19.12.22'
$ws.Cells.Item(72, 1).Value = 'Add Document Types'
$ws.Cells.Item(72, 2).Value = 'FAILED'
$ws.Cells.Item(72, 3).Value = 'chrome'
$ws.Cells.Item(72, 4).Value = '19.12.22'
$ws.Cells.Item(73, 1).Value = 'Edit Document Types'
$ws.Cells.Item(73, 2).Value = 'FAILED'
$ws.Cells.Item(73, 3).Value = 'chrome'
$ws.Cells.Item(73, 4).Value = '19.12.22'
$ws.Cells.Item(74, 1).Value = 'Add Document Types'
$ws.Cells.Item(74, 2).Value = 'PASSED'
$ws.Cells.Item(74, 3).Value = 'chrome'
$ws.Cells.Item(74, 4).Value = '20.12.22'
$ws.Cells.Item(75, 1).Value = 'Add Document Types'
$ws.Cells.Item(75, 2).Value = 'PASSED'
$ws.Cells.Item(75, 3).Value = 'chrome'
$ws.Cells.Item(75, 4).Value = '20.12.22'
$ws.Cells.Item(76, 1).Value = 'Edit Document Types'
$ws.Cells.Item(76, 2).Value = 'FAILED'
$ws.Cells.Item(76, 3).Value = 'chrome'
$ws.Cells.Item(76, 4).Value = '20.12.22'
$ws.Cells.Item(77, 1).Value = 'Delete Document Types'
$ws.Cells.Item(77, 2).Value = 'PASSED'
$ws.Cells.Item(77, 3).Value = 'chrome'
$ws.Cells.Item(77, 4).Value = '20.12.22'
$ws.Cells.Item(78, 1).Value = 'Negative Delete Document Types'
$ws.Cells.Item(78, 2).Value = 'PASSED'
$ws.Cells.Item(78, 3).Value = 'chrome'
$ws.Cells.Item(78, 4).Value = '20.12.22'
$ws.Cells.Item(79, 1).Value = 'Login with valid username and password'
$ws.Cells.Item(79, 2).Value = 'PASSED'
$ws.Cells.Item(79, 3).Value = 'chrome'
$ws.Cells.Item(79, 4).Value = '22.12.22'
$ws.Cells.Item(80, 1).Value = 'Login with valid username and password'
$ws.Cells.Item(80, 2).Value = 'PASSED'
$ws.Cells.Item(80, 3).Value = 'chrome'
$ws.Cells.Item(80, 4).Value = '22.12.22'
$ws.Cells.Item(81, 1).Value = 'Add Document Types'
$ws.Cells.Item(81, 2).Value = 'PASSED'
$ws.Cells.Item(81, 3).Value = 'chrome'
$ws.Cells.Item(81, 4).Value = '17.01.23'
$ws.Cells.Item(82, 1).Value = 'Edit Document Types'
$ws.Cells.Item(82, 2).Value = 'PASSED'
$ws.Cells.Item(82, 3).Value = 'chrome'
$ws.Cells.Item(82, 4).Value = '17.01.23'
$ws.Cells.Item(83, 1).Value = 'Delete Document Types'
$ws.Cells.Item(83, 2).Value = 'PASSED'
$ws.Cells.Item(83, 3).Value = 'chrome'
$ws.Cells.Item(83, 4).Value = '17.01.23'
$ws.Cells.Item(84, 1).Value = 'Create A Citizenship'
$ws.Cells.Item(84, 2).Value = 'PASSED'
$ws.Cells.Item(84, 3).Value = 'chrome'
$ws.Cells.Item(84, 4).Value = '17.01.23'
$ws.Cells.Item(85, 1).Value = 'Create An Existant Citizenship'
$ws.Cells.Item(85, 2).Value = 'PASSED'
$ws.Cells.Item(85, 3).Value = 'chrome'
$ws.Cells.Item(85, 4).Value = '17.01.23'
$ws.Cells.Item(86, 1).Value = 'Update the Citizenship'
$ws.Cells.Item(86, 2).Value = 'PASSED'
$ws.Cells.Item(86, 3).Value = 'chrome'
$ws.Cells.Item(86, 4).Value = '17.01.23'
$ws.Cells.Item(87, 1).Value = 'Delete the Citizenship'
$ws.Cells.Item(87, 2).Value = 'PASSED'
$ws.Cells.Item(87, 3).Value = 'chrome'
$ws.Cells.Item(87, 4).Value = '17.01.23'
$ws.Cells.Item(88, 1).Value = 'Search and delete an unavailable Citizenship'
$ws.Cells.Item(88, 2).Value = 'PASSED'
$ws.Cells.Item(88, 3).Value = 'chrome'
$ws.Cells.Item(88, 4).Value = '17.01.23'
$ws.Cells.Item(89, 1).Value = 'Create Country'
$ws.Cells.Item(89, 2).Value = 'PASSED'
$ws.Cells.Item(89, 3).Value = 'chrome'
$ws.Cells.Item(89, 4).Value = '17.01.23'
$ws.Cells.Item(90, 1).Value = 'Create a Attestations'
$ws.Cells.Item(90, 2).Value = 'FAILED'
$ws.Cells.Item(90, 3).Value = 'chrome'
$ws.Cells.Item(90, 4).Value = '17.01.23'
$ws.Cells.Item(91, 1).Value = 'Edit a Attestations'
$ws.Cells.Item(91, 2).Value = 'FAILED'
$ws.Cells.Item(91, 3).Value = 'chrome'
$ws.Cells.Item(91, 4).Value = '17.01.23'
$ws.Cells.Item(92, 1).Value = 'Delete a Attestations'
$ws.Cells.Item(92, 2).Value = 'FAILED'
$ws.Cells.Item(92, 3).Value = 'chrome'
$ws.Cells.Item(92, 4).Value = '17.01.23'
$ws.Cells.Item(93, 1).Value = 'Add Document Types'
$ws.Cells.Item(93, 2).Value = 'PASSED'
$ws.Cells.Item(93, 3).Value = 'chrome'
$ws.Cells.Item(93, 4).Value = '17.01.23'
$ws.Cells.Item(94, 1).Value = 'Add Document Types'
$ws.Cells.Item(94, 2).Value = 'FAILED'
$ws.Cells.Item(94, 3).Value = 'chrome'
$ws.Cells.Item(94, 4).Value = '17.01.23'
$ws.Cells.Item(95, 1).Value = 'Add Document Types'
$ws.Cells.Item(95, 2).Value = 'FAILED'
$ws.Cells.Item(95, 3).Value = 'chrome'
$ws.Cells.Item(95, 4).Value = '17.01.23'
$ws.Cells.Item(96, 1).Value = 'Add Document Types'
$ws.Cells.Item(96, 2).Value = 'PASSED'
$ws.Cells.Item(96, 3).Value = 'chrome'
$ws.Cells.Item(96, 4).Value = '17.01.23'
$ws.Cells.Item(97, 1).Value = 'Add Document Types'
$ws.Cells.Item(97, 2).Value = 'PASSED'
$ws.Cells.Item(97, 3).Value = 'chrome'
$ws.Cells.Item(97, 4).Value = '24.01.23'
$ws.Cells.Item(98, 1).Value = 'Add School Locations'
$ws.Cells.Item(98, 2).Value = 'PASSED'
$ws.Cells.Item(98, 3).Value = 'chrome'
$ws.Cells.Item(98, 4).Value = '24.01.23'
$ws.Cells.Item(99, 1).Value = 'Add School Locations'
$ws.Cells.Item(99, 2).Value = 'PASSED'
$ws.Cells.Item(99, 3).Value = 'chrome'
$ws.Cells.Item(99, 4).Value = '24.01.23'
$ws.Cells.Item(100, 1).Value = 'Add School Locations'
$ws.Cells.Item(100, 2).Value = 'PASSED'
$ws.Cells.Item(100, 3).Value = 'chrome'
$ws.Cells.Item(100, 4).Value = '24.01.23'
$ws.Cells.Item(101, 1).Value = 'Add School Locations'
$ws.Cells.Item(101, 2).Value = 'PASSED'
$ws.Cells.Item(101, 3).Value = 'chrome'
$ws.Cells.Item(101, 4).Value = '24.01.23'
